# Populate the newly-added "nearest_upstream_orf" (M) / "nearest_downstream_orf" (N)
# flanking-gene data for the proto-refseqs rows on Sheet1. Previously these cells
# held the placeholder "NK"; this fills in the actual flanking ORF/gene identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = "AABR07048228.1"
$ws.Range("N2").Value = "AABR07048231.1"

$ws.Range("M3").Value = "PDZRN4"
$ws.Range("N3").Value = "PDZRN4"

$ws.Range("M4").Value = "GM14461"
$ws.Range("N4").Value = "UBE2E3"

$ws.Range("M5").Value = "CCDC146"
$ws.Range("N5").Value = "CCDC146"

$ws.Range("M7").Value = "U6"
$ws.Range("N7").Value = "ENSMODG00000006365"

$ws.Range("M9").Value = "ENSMODG00000004341"
$ws.Range("N9").Value = "ENSMODG00000045240"

$ws.Range("M10").Value = "U6"
$ws.Range("N10").Value = "ENSMODG00000006365"

$ws.Range("N15").Value = "ENSCPOG00000022174"

$ws.Range("M17").Value = "ENSCPOG00000033273"
$ws.Range("N17").Value = "ENSCPOG00000026716"

$ws.Range("M40").Value = "KLHL1"
$ws.Range("N40").Value = "KLHL1"

$ws.Range("M41").Value = "ENSSHAG00000023532"
$ws.Range("N41").Value = "ENSSHAG00000027488"

$ws.Range("M47").Value = "ME3"
$ws.Range("N47").Value = "ME3"

$ws.Range("M52").Value = "AGPAT5"

$ws.Range("M53").Value = "U6"
$ws.Range("N53").Value = "ENSPCIG0000036639"

$ws.Range("M55").Value = "ENSVURG00010023687"
$ws.Range("N55").Value = "FZD8"

$ws.Range("M59").Value = "ENSMODG00000007612"
$ws.Range("N59").Value = "ENSMODG00000047647"

$ws.Range("M60").Value = "ENSMODG00000035635"
$ws.Range("N60").Value = "ENSMODG00000042352"

$ws.Range("M67").Value = "ENPP2"
$ws.Range("N67").Value = "ENPP2"

$ws.Range("M68").Value = "KLF6"
$ws.Range("N68").Value = "ENSPCIG00000033363"

$ws.Range("M69").Value = "BTBD1"
$ws.Range("N69").Value = "BTBD1"

$ws.Range("M70").Value = "VPS8"
$ws.Range("N70").Value = "VPS8"

$ws.Range("M71").Value = "AGTPBP1"
$ws.Range("N71").Value = "AGTPBP1"

$ws.Range("M72").Value = "ENSPCIG00000032522"
$ws.Range("N72").Value = "ENSPCIG00000024782"

$ws.Range("M73").Value = "ENSPCIG00000015775"
$ws.Range("N73").Value = "ITGA6"

$ws.Range("M74").Value = "SFRP2"
$ws.Range("N74").Value = "ENSPCIG00000029650"

$ws.Range("M75").Value = "ENSPCIG00000031494"
$ws.Range("N75").Value = "ENSPCIG00000011842"

$ws.Range("M81").Value = "ENSSHAG00000017845"
$ws.Range("N81").Value = "PCDH18"

$ws.Range("M82").Value = "ENSSHAG00000026805"
$ws.Range("N82").Value = "ENSSHAG00000020835"

$ws.Range("M83").Value = "ENSSHAG00000002243"
$ws.Range("N83").Value = "ENSSHAG00000028478"

$ws.Range("M85").Value = "ENSMODG00000036286"
$ws.Range("N85").Value = "ENSMODG00000036286"

$ws.Range("M96").Value = "VPS8"
$ws.Range("N96").Value = "VPS8"

$ws.Range("M104").Value = "ENSVURG00010005697"
$ws.Range("N104").Value = "NT5E"

$ws.Range("M105").Value = "ENPP2"
$ws.Range("N105").Value = "ENPP2"

# Widen the M:N columns to fit the new longer identifiers (~24.16 chars), matching
# the column-width metadata added alongside this data in the source workbook.
$ws.Range("M1:N1").ColumnWidth = 24.1640625

# Leave the selection where the author last left it after entering this data.
$ws.Range("M106").Select() | Out-Null

Write-Output "Flanking ORF data populated for proto-refseqs rows."
